$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.145.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.247.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '395.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.590'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.243.46'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.625'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '39.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0981'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.47%  '

$ws.Range("E13").Value = '  +1.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.765.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.249.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.04%  '

$ws.Range("E18").Value = '  -2.89%  '

$ws.Range("E19").Value = '  +2.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '57.071.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000112'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '295.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.77%  '

$ws.Range("E30").Value = '  -1.70%  '

$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.91%  '

$ws.Range("E33").Value = '  -4.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '39.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.26%  '

$ws.Range("E35").Value = '  -3.61%  '

$ws.Range("E36").Value = '  +2.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.27%  '

$ws.Range("E40").Value = '  +1.88%  '

$ws.Range("E41").Value = '  +4.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '134.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.280'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.16%  '

$ws.Range("E48").Value = '  +3.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.158.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("E50").Value = '  -4.34%  '

$ws.Range("E51").Value = '  +15.81%  '
